$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.022.65'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.911.82'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('D4').Value = "'0.9990"
$ws.Range('E4').Value = '  -0.68%  '
$ws.Range('D5').Value = "'324.11"
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = "'0.3870"
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = "'0.07826"
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').Value = "'0.9895"
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').Value = "'21.90"
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').Value = '1.889.23'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = "'5.765"
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = "'7.009"
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').Value = "'0.07055"
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = "'87.42"
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = "'1.002"
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').Value = "'0.000009904"
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').Value = "'17.04"
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').Value = '29.046.99'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = "'5.381"
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').Value = "'11.10"
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').Value = '2.141.57'
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').Value = "'156.09"
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = "'19.31"
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').Value = "'5.847"
$ws.Range('E28').Value = '  -2.66%  '
$ws.Range('D29').Value = "'118.23"
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = "'1.852"
$ws.Range('E30').Value = '  -4.40%  '
$ws.Range('D31').Value = "'0.09304"
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').Value = "'0.8797"
$ws.Range('E32').Value = '  -3.67%  '
$ws.Range('D33').Value = "'5.184"
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('E34').Value = '  -2.54%  '
$ws.Range('D35').Value = "'3.122"
$ws.Range('E35').Value = '  -4.96%  '
$ws.Range('D36').Value = "'0.05759"
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').Value = "'1.168"
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').Value = "'0.02088"
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').Value = "'0.9996"
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').Value = "'0.5680"
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').Value = "'7.645"
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').Value = "'0.1807"
$ws.Range('E42').Value = '  +1.27%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = "'0.000002986"
$ws.Range('E43').Value = '  +84.23%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'9.663"
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = "'2.205"
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = "'0.5316"
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = "'0.06924"
$ws.Range('D49').Value = "'1.833"
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('D50').Value = "'2.558"
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = "'112.52"
$ws.Range('E51').Value = '  -0.33%  '

$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
